$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.320.93"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +1.93%  '

$ws.Cells.Item(3, 4).Value = "'1.649.75"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.57%  '

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.15%  '

$ws.Cells.Item(5, 4).Value = "'217.54"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.62%  '

$ws.Cells.Item(6, 4).Value = "'0.508"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.35%  '

$ws.Cells.Item(7, 5).Value = '  -0.17%  '

$ws.Cells.Item(8, 5).Value = '  +0.02%  '

$ws.Cells.Item(9, 5).Value = '  +0.30%  '

$ws.Cells.Item(10, 4).Value = "'20.00"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.33%  '

$ws.Cells.Item(11, 5).Value = '  +0.13%  '

$ws.Cells.Item(12, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(12, 4).Value = "'1.878.94"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.65%  '

$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).Value = "'4.31"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.84%  '

$ws.Cells.Item(14, 4).Value = "'1.674.91"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.03%  '

$ws.Cells.Item(15, 4).Value = "'0.551"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.09%  '

$ws.Cells.Item(16, 5).Value = '  +0.09%  '

$ws.Cells.Item(17, 4).Value = "'63.62"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.87%  '

$ws.Cells.Item(18, 4).Value = "'26.314.52"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +1.76%  '

$ws.Cells.Item(19, 5).Value = '  -0.14%  '

$ws.Cells.Item(20, 4).Value = "'196.95"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +2.11%  '

$ws.Cells.Item(21, 5).Value = '  -0.38%  '

$ws.Cells.Item(22, 5).Value = '  +1.11%  '

$ws.Cells.Item(23, 5).Value = '  -0.14%  '

$ws.Cells.Item(24, 5).Value = '  -2.07%  '

$ws.Cells.Item(25, 4).Value = "'143.10"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.49%  '

$ws.Cells.Item(27, 5).Value = '  +1.43%  '

$ws.Cells.Item(28, 5).Value = '  +0.77%  '

$ws.Cells.Item(29, 5).Value = '  +0.82%  '

$ws.Cells.Item(30, 4).Value = "'1.25"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.14%  '

$ws.Cells.Item(31, 4).Value = "'0.0507"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.54%  '

$ws.Cells.Item(32, 5).Value = '  +0.51%  '

$ws.Cells.Item(34, 5).Value = '  +2.36%  '

$ws.Cells.Item(35, 5).Value = '  +1.03%  '

$ws.Cells.Item(36, 4).Value = "'0.918"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.15%  '

$ws.Cells.Item(37, 4).Value = "'0.557"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.89%  '

$ws.Cells.Item(38, 4).Value = "'1.137.36"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.17%  '

$ws.Cells.Item(39, 5).Value = '  -1.41%  '

$ws.Cells.Item(40, 5).Value = '  +1.02%  '

$ws.Cells.Item(41, 5).Value = '  -0.22%  '

$ws.Cells.Item(42, 4).Value = "'5.67"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.30%  '

$ws.Cells.Item(43, 4).Value = "'100.53"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.25%  '

$ws.Cells.Item(44, 4).Value = "'0.803"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.56%  '

$ws.Cells.Item(45, 4).Value = "'1.788.00"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.68%  '

$ws.Cells.Item(46, 4).Value = "'56.51"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.98%  '

$ws.Cells.Item(47, 5).Value = '  +3.73%  '

$ws.Cells.Item(48, 5).Value = '  +3.18%  '

$ws.Cells.Item(49, 4).Value = "'7.72"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.07%  '

$ws.Cells.Item(50, 5).Value = '  -0.02%  '

$ws.Cells.Item(51, 4).Value = "'0.0975"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.33%  '
